$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title text updated from Bloomberg BCS branding to Ceph
$ws.Range("A1").Value = "Ceph Cloud Storage Erasure Coding Calculations"

# Input values updated for the Jewel upgrade sizing
$ws.Range("B2").Value = 1224
$ws.Range("B5").Value = 8

# Move active selection to B3, matching the saved selection state
$ws.Range("B3").Select()
